$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for columns G:T, rows 2-9 (14 columns x 8 rows)
$newValues = @(
        @(0.3379886666666667,1.013966,0.4664623486936776,0.4664623486936776,1,0.3333333333333333,0.1199136666666667,0.359741,0.03475238343556118,0.03475238343556118,0.04052946031177777,0.364765142806,0.01621067840005512,0.01621067840005512),
        @(0.3379886666666667,1.013966,0.4664623486936776,0.4664623486936776,2,0.6666666666666666,0.247271,0.7418130000000001,0.07166202855244176,0.07166202855244176,0.08357479559533333,0.752173160358,0.03342763815072537,0.03342763815072537),
        @(0.3379886666666667,1.013966,0.4664623486936776,0.4664623486936776,3,1,0.355759,1.067277,0.1031031201224087,0.1031031201224087,0.1202425100646667,1.082182590582,0.04809372356994514,0.04809372356994514),
        @(0.3379886666666667,1.013966,0.4664623486936776,0.4664623486936776,3,1,2.727572666666667,8.182718000000001,0.7904824678895883,0.7904824678895883,0.9218886488431113,8.296997839588,0.3687303085729519,0.3687303085729519),
        @(0.38659,1.15977,0.5335376513063224,0.5335376513063224,1,0.3333333333333333,0.1199136666666667,0.359741,0.03475238343556118,0.03475238343556118,0.04635742439666666,0.41721681957,0.01854170503550606,0.01854170503550606),
        @(0.38659,1.15977,0.5335376513063224,0.5335376513063224,2,0.6666666666666666,0.247271,0.7418130000000001,0.07166202855244176,0.07166202855244176,0.09559249589,0.86033246301,0.03823439040171639,0.03823439040171639),
        @(0.38659,1.15977,0.5335376513063224,0.5335376513063224,3,1,0.355759,1.067277,0.1031031201224087,0.1031031201224087,0.13753287181,1.23779584629,0.05500939655246357,0.05500939655246357),
        @(0.38659,1.15977,0.5335376513063224,0.5335376513063224,3,1,2.727572666666667,8.182718000000001,0.7904824678895883,0.7904824678895883,1.054452317206667,9.49007085486,0.4217521593166363,0.4217521593166363)
    )

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 7  # column G is index 7
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
